$wb = $excel.ActiveWorkbook

# Data describing the cell updates per worksheet: row number -> column letter -> new value
# A value of $null means the cell should be cleared (no stored value).
$updates = @{}

$updates["ALC"] = @{
    16 = @{ "H" = 1833.2222; "I" = 999; "J" = 1937.5; "K" = 999; "L" = 1937.5; "M" = -769; "N" = -2397.5 }
    49 = @{ "H" = 2955.5557; "J" = 4066.6667; "L" = 12200.0001; "N" = -12472.0001 }
    51 = @{ "H" = 76940290; "I" = 333334850; "K" = 333334850; "M" = -333334366 }
    74 = @{ "H" = 3521.739; "I" = 3521.739; "J" = 0; "K" = 3521.739; "L" = 0; "M" = -2585.739; "N" = $null }
    76 = @{ "H" = 18870884; "I" = 21742068; "J" = 3100; "K" = 21742068; "L" = 3100; "M" = -21741753; "N" = -3730 }
    77 = @{ "H" = 3521.739; "I" = 3521.739; "J" = 0; "K" = 17608.695; "L" = 0; "M" = -12928.695; "N" = $null }
    79 = @{ "H" = 18870884; "I" = 21742068; "J" = 3100; "K" = 21742068; "L" = 3100; "M" = -21740976; "N" = -5284 }
}

$updates["ARM"] = @{
    32 = @{ "H" = 7252616.5; "I" = 6380.5356; "K" = 6380.5356; "M" = -6093.5356 }
    58 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
    62 = @{ "H" = 25000; "J" = 25000; "L" = 25000; "N" = -26372 }
    63 = @{ "H" = 3746.6667; "I" = 2500; "J" = 4993.3335; "K" = 2500; "L" = 4993.3335; "M" = -1814; "N" = -6365.3335 }
    65 = @{ "H" = 25000; "J" = 25000; "L" = 75000; "N" = -81864 }
    66 = @{ "H" = 3746.6667; "I" = 2500; "J" = 4993.3335; "K" = 12500; "L" = 24966.6675; "M" = -9068; "N" = -31830.6675 }
    70 = @{ "H" = 39800; "J" = 39800; "L" = 39800; "N" = -40430 }
    73 = @{ "H" = 39800; "J" = 39800; "L" = 39800; "N" = -41984 }
    76 = @{ "H" = 24500; "J" = 24500; "L" = 24500; "N" = -25176 }
    79 = @{ "H" = 24500; "J" = 24500; "L" = 24500; "N" = -26840 }
}

$updates["BSM"] = @{
    76 = @{ "H" = 10314; "J" = 10314; "L" = 10314; "N" = -10944 }
    79 = @{ "H" = 10314; "J" = 10314; "L" = 10314; "N" = -12498 }
    88 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
    91 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
}

$updates["CRP"] = @{
    32 = @{ "H" = 1000; "I" = 0; "J" = 1000; "K" = 0; "L" = 1000; "M" = $null; "N" = -1632 }
    82 = @{ "H" = 25499.75; "J" = 25499.75; "L" = 25499.75; "N" = -26221.75 }
    85 = @{ "H" = 25499.75; "J" = 25499.75; "L" = 25499.75; "N" = -27995.75 }
    88 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
    91 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
}

$updates["CUL"] = @{
    26 = @{ "H" = 163.64285; "I" = 66.454544; "J" = 520; "K" = 199.363632; "L" = 1560; "M" = 88.636368; "N" = -2136 }
}

$updates["GSM"] = @{
    82 = @{ "H" = 10150.5; "J" = 10150.5; "L" = 10150.5; "N" = -10916.5 }
    85 = @{ "H" = 10150.5; "J" = 10150.5; "L" = 10150.5; "N" = -12802.5 }
    88 = @{ "H" = 17999; "J" = 17999; "L" = 17999; "N" = -18855 }
    91 = @{ "H" = 17999; "J" = 17999; "L" = 17999; "N" = -20963 }
}

$updates["LTW"] = @{
    47 = @{ "H" = 13191.667; "J" = 13191.667; "L" = 13191.667; "N" = -14171.667 }
    52 = @{ "H" = 13191.667; "J" = 13191.667; "L" = 13191.667; "N" = -13657.667 }
    64 = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null; "N" = $null }
    67 = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null; "N" = $null }
    76 = @{ "H" = 23180.125; "J" = 23180.125; "L" = 23180.125; "N" = -23856.125 }
    79 = @{ "H" = 23180.125; "J" = 23180.125; "L" = 23180.125; "N" = -25520.125 }
    82 = @{ "H" = 4360.778; "I" = 1676; "J" = 6069.273; "K" = 1676; "L" = 6069.273; "M" = -1315; "N" = -6791.273 }
    85 = @{ "H" = 4360.778; "I" = 1676; "J" = 6069.273; "K" = 1676; "L" = 6069.273; "M" = -428; "N" = -8565.273000000001 }
}

$updates["WVR"] = @{
    76 = @{ "H" = 19324.334; "J" = 19324.334; "L" = 19324.334; "N" = -19954.334 }
    79 = @{ "H" = 19324.334; "J" = 19324.334; "L" = 19324.334; "N" = -21508.334 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($rowNum in $rowsForSheet.Keys) {
        $colsForRow = $rowsForSheet[$rowNum]
        foreach ($colLetter in $colsForRow.Keys) {
            $cellAddress = "$colLetter$rowNum"
            $newValue = $colsForRow[$colLetter]
            if ($null -eq $newValue) {
                $ws.Range($cellAddress).ClearContents()
            } else {
                $ws.Range($cellAddress).Value = $newValue
            }
        }
    }
}

Write-Host "Applied $($updates.Count) sheet updates."